$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the coordinate table (rows 31-41, columns B/Nr., C/Y, D/X, I/Sipërfaqja) ---

# Row 31
$ws.Range("B31").Value = 50
$ws.Range("C31").Value = 7510648.7396
$ws.Range("D31").Value = 4693483.1380000003
$ws.Range("I31").Value = 79.843999999999994

# Row 32
$ws.Range("B32").Value = 51
$ws.Range("C32").Value = 7510648.5892000003
$ws.Range("D32").Value = 4693483.0460999999

# Row 33
$ws.Range("B33").Value = 52
$ws.Range("C33").Value = 7510641.7006000001
$ws.Range("D33").Value = 4693478.8332000002

# Row 34
$ws.Range("B34").Value = 60
$ws.Range("C34").Value = 7510653.6935999999
$ws.Range("D34").Value = 4693474.8667000001

# Row 35
$ws.Range("B35").Value = 61
$ws.Range("C35").Value = 7510649.8039999995
$ws.Range("D35").Value = 4693472.3590000002

# Row 36
$ws.Range("B36").Value = 62
$ws.Range("C36").Value = 7510648.0659999996
$ws.Range("D36").Value = 4693471.3150000004

# Row 37
$ws.Range("B37").Value = 63
$ws.Range("C37").Value = 7510648.0179000003
$ws.Range("D37").Value = 4693471.3956000004

# Row 38
$ws.Range("B38").Value = 64
$ws.Range("C38").Value = 7510646.6425000001
$ws.Range("D38").Value = 4693470.5544999996

# Row 39
$ws.Range("B39").Value = 75
$ws.Range("C39").Value = 7510649.7476000004
$ws.Range("D39").Value = 4693472.4534999998
$ws.Range("I39").Value = 8.3149999999999995

# Row 40
$ws.Range("B40").Value = 43
$ws.Range("C40").Value = 7510641.1344999997
$ws.Range("D40").Value = 4693479.7813999997

# Row 41
$ws.Range("B41").Value = 44
$ws.Range("C41").Value = 7510648.0900999997
$ws.Range("D41").Value = 4693483.8622000003

# --- Update the sheet view / selection ---
$app = $ws.Application
$win = $app.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("B30:I41").Select()
